$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "address.lane"
$ws.Range("E2").Value = "lane 1"
$ws.Range("E3").Value = "lane 2"

$ws.Range("F1").Value = "contact.phone"
$ws.Range("G1").Value = "contact.email"

$ws.Range("F2").Value = 9206918946
$ws.Range("G2").Value = "deepak.kumar@gmail.com"

$ws.Range("F3").Value = 9206918947
$ws.Range("G3").Value = "kumar.deepak@gmail.com"

$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(6).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(7).EntireColumn.AutoFit() | Out-Null

$ws.Range("G3").Select() | Out-Null
